$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 42,4

$data[0,0] = 'DESN'
$data[0,1] = 2
$data[0,2] = 1
$data[0,3] = 'Follows standard approaches and established design patterns to create new designs for simple systems or system components'
$data[1,0] = 'DESN'
$data[1,1] = 2
$data[1,2] = 2
$data[1,3] = 'Identifies and resolves minor design issues'
$data[2,0] = 'DESN'
$data[2,1] = 2
$data[2,2] = 3
$data[2,3] = 'Identifies alternative design options and seeks guidance when deviating from established design patterns'
$data[3,0] = 'SWDN'
$data[3,1] = 2
$data[3,2] = 1
$data[3,3] = 'Undertakes complete design of moderately complex software applications or components'
$data[4,0] = 'SWDN'
$data[4,1] = 2
$data[4,2] = 2
$data[4,3] = 'Applies agreed standards, guidelines, patterns and tools'
$data[5,0] = 'SWDN'
$data[5,1] = 2
$data[5,2] = 3
$data[5,3] = 'Assists as part of a team in the design of components of larger software systems'
$data[6,0] = 'SWDN'
$data[6,1] = 2
$data[6,2] = 4
$data[6,3] = 'Specifies user and/or system interfaces'
$data[7,0] = 'SWDN'
$data[7,1] = 2
$data[7,2] = 5
$data[7,3] = 'Creates multiple design views to address the different stakeholders'' concerns and to handle functional and non-functional requirements'
$data[8,0] = 'SWDN'
$data[8,1] = 2
$data[8,2] = 6
$data[8,3] = 'Assists in the evaluation of options and trade-offs'
$data[9,0] = 'SWDN'
$data[9,1] = 2
$data[9,2] = 7
$data[9,3] = 'Collaborates in reviews of work with others as appropriate'
$data[10,0] = 'PROG'
$data[10,1] = 2
$data[10,2] = 1
$data[10,3] = 'Designs, codes, verifies, tests, documents, amends and refactors moderately complex programs/scripts'
$data[11,0] = 'PROG'
$data[11,1] = 2
$data[11,2] = 2
$data[11,3] = 'Applies agreed standards and tools to achieve a well-engineered result'
$data[12,0] = 'PROG'
$data[12,1] = 2
$data[12,2] = 3
$data[12,3] = 'Monitors and reports on progress'
$data[13,0] = 'PROG'
$data[13,1] = 2
$data[13,2] = 4
$data[13,3] = 'Identifies issues related to software development activities'
$data[14,0] = 'PROG'
$data[14,1] = 2
$data[14,2] = 5
$data[14,3] = 'Proposes practical solutions to resolve issues'
$data[15,0] = 'PROG'
$data[15,1] = 2
$data[15,2] = 6
$data[15,3] = 'Collaborates in reviews of work with others as appropriate'
$data[16,0] = 'SINT'
$data[16,1] = 2
$data[16,2] = 1
$data[16,3] = 'Defines the software modules needed for an integration build and produces a build definition for each generation of the software'
$data[17,0] = 'SINT'
$data[17,1] = 2
$data[17,2] = 2
$data[17,3] = 'Accepts completed software modules, ensuring that they meet defined criteria'
$data[18,0] = 'SINT'
$data[18,1] = 2
$data[18,2] = 3
$data[18,3] = 'Produces software builds from software source code for loading onto target hardware'
$data[19,0] = 'SINT'
$data[19,1] = 2
$data[19,2] = 4
$data[19,3] = 'Configures the hardware and software environment as required by the system being integrated'
$data[20,0] = 'SINT'
$data[20,1] = 2
$data[20,2] = 5
$data[20,3] = 'Produces integration test specifications, conducts tests and records and reports on outcomes'
$data[21,0] = 'SINT'
$data[21,1] = 2
$data[21,2] = 6
$data[21,3] = 'Diagnoses faults and records and reports on the results of tests'
$data[22,0] = 'SINT'
$data[22,1] = 2
$data[22,2] = 7
$data[22,3] = 'Produces system integration reports'
$data[23,0] = 'TEST'
$data[23,1] = 2
$data[23,2] = 1
$data[23,3] = 'Designs test cases and test scripts under own direction, mapping back to pre-determined criteria, recording and reporting test outcomes'
$data[24,0] = 'TEST'
$data[24,1] = 2
$data[24,2] = 2
$data[24,3] = 'Participates in requirement, design and specification reviews, and uses this information to design test plans and test conditions'
$data[25,0] = 'TEST'
$data[25,1] = 2
$data[25,2] = 3
$data[25,3] = 'Applies agreed standards to specify and perform manual and automated testing'
$data[26,0] = 'TEST'
$data[26,1] = 2
$data[26,2] = 4
$data[26,3] = 'Automates testing tasks and builds test coverage through existing or new infrastructure'
$data[27,0] = 'TEST'
$data[27,1] = 2
$data[27,2] = 5
$data[27,3] = 'Analyses and reports on test activities, results, issues and risks'
$data[28,0] = 'CFMG'
$data[28,1] = 2
$data[28,2] = 1
$data[28,3] = 'Applies tools, techniques and processes to track, log and correct information related to configuration items'
$data[29,0] = 'CFMG'
$data[29,1] = 2
$data[29,2] = 2
$data[29,3] = 'Verifies and approves changes ensuring the protection of assets and components from unauthorised change, diversion and inappropriate use'
$data[30,0] = 'CFMG'
$data[30,1] = 2
$data[30,2] = 3
$data[30,3] = 'Ensures that users comply with identification standards for object types, environments, processes, life cycles, documentation, versions, formats, baselines, releases and templates'
$data[31,0] = 'CFMG'
$data[31,1] = 2
$data[31,2] = 4
$data[31,3] = 'Performs audits to check the accuracy of the information and undertakes any necessary corrective action under direction'
# row 47 intentionally left blank
$data[32,0] = $null
$data[32,1] = $null
$data[32,2] = $null
$data[32,3] = $null
$data[33,0] = 'REQM'
$data[33,1] = 2
$data[33,2] = 1
$data[33,3] = 'Defines and manages scoping, requirements definition and prioritisation activities for small-scale changes and assists with more complex change initiatives'
$data[34,0] = 'REQM'
$data[34,1] = 2
$data[34,2] = 2
$data[34,3] = 'Follows agreed standards and applies appropriate techniques to elicit and document detailed requirements'
$data[35,0] = 'REQM'
$data[35,1] = 2
$data[35,2] = 3
$data[35,3] = 'Provides constructive challenge to stakeholders as required'
$data[36,0] = 'REQM'
$data[36,1] = 2
$data[36,2] = 4
$data[36,3] = 'Reviews requirements for errors and omissions'
$data[37,0] = 'REQM'
$data[37,1] = 2
$data[37,2] = 5
$data[37,3] = 'Prioritises requirements and documents traceability to source'
$data[38,0] = 'REQM'
$data[38,1] = 2
$data[38,2] = 6
$data[38,3] = 'Provides input to the requirements base-line'
$data[39,0] = 'REQM'
$data[39,1] = 2
$data[39,2] = 7
$data[39,3] = 'Investigates, manages and applies authorised requests for changes to base-lined requirements, in line with change management policy'
$data[40,0] = 'RELM'
$data[40,1] = 2
$data[40,2] = 1
$data[40,3] = 'Uses approved tools and techniques for specific deployment activities'
$data[41,0] = 'RELM'
$data[41,1] = 2
$data[41,2] = 2
$data[41,3] = 'Administers the recording of activities, logging of results and documents technical activities undertaken'

$ws.Range("A15:D56").Value2 = $data

Write-Host "Done writing rows 15 to 56"